$d = $word.ActiveDocument

# Find the existing "Variations" bullet that reads
# "Having 13 cards and more Contracts" so the new bullet can be appended
# right after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Having 13 cards and more Contracts",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $anchorPara = $rng.Paragraphs(1)

    # Insert a brand-new paragraph right after the anchor. Word copies the
    # anchor's paragraph formatting (ListParagraph style, numId 6 / ilvl 0
    # bullet numbering) onto the freshly inserted paragraph automatically.
    $anchorPara.Range.InsertParagraphAfter()

    # The paragraph we just inserted is now the document's last paragraph
    # (the anchor was the last bullet in the "Variations" list, right
    # before the closing sectPr). Fill in its text.
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = "Allow runs to be mixed suit"
}
